$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Prepend the generator's output folder and append the .pdf extension to
# every file-name cell in column A (rows 2-6), turning the bare file
# names into full Windows paths to the generated PDF.
$prefix = "C:\Users\Lucas.Aguiar\Desktop\gerador_senha\pdf\"

for ($r = 2; $r -le 6; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $fileName = $cell.Value2
    $cell.Value = "$prefix$fileName.pdf"
}

# Column A now holds much longer text, so widen it to fit ("best fit"
# style width); column B is left untouched.
$ws.Columns.Item(1).ColumnWidth = 65.14

# Move the active selection, matching the author's last cursor position.
$ws.Range("D14").Select()
